$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Target cluster changes from "ECs" to "FAPs"
$ws.Range("D2").Value = "FAPs"
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.029119
$ws.Range("N2").Value = 0.087357
$ws.Range("O2").Value = 0.4059019501247578
$ws.Range("P2").Value = 0.4059019501247578
$ws.Range("Q2").Value = 0.004192786572000001
$ws.Range("R2").Value = 0.03773507914800001
$ws.Range("S2").Value = 0.4059019501247578
$ws.Range("T2").Value = 0.4059019501247578

# Row 3: Target cluster changes from "FAPs" to "MuSCs"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.04229766666666667
$ws.Range("N3").Value = 0.126893
$ws.Range("O3").Value = 0.5896049103927664
$ws.Range("P3").Value = 0.5896049103927664
$ws.Range("Q3").Value = 0.006090356428000001
$ws.Range("R3").Value = 0.054813207852
$ws.Range("S3").Value = 0.5896049103927664
$ws.Range("T3").Value = 0.5896049103927664

# Row 4: Target cluster changes from "MuSCs" to "Neutrophils"
$ws.Range("D4").Value = "Neutrophils"
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.0003223333333333333
$ws.Range("N4").Value = 0.000967
$ws.Range("O4").Value = 0.004493139482475827
$ws.Range("P4").Value = 0.004493139482475827
$ws.Range("Q4").Value = 0.000046412132
$ws.Range("R4").Value = 0.000417709188
$ws.Range("S4").Value = 0.004493139482475827
$ws.Range("T4").Value = 0.004493139482475827
